$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 173.90323
$ws.Range("I9").Value = 72.44444
$ws.Range("K9").Value = 72.44444
$ws.Range("M9").Value = 96.55556
$ws.Range("H19").Value = 853.4545000000001
$ws.Range("I19").Value = 1085.8572
$ws.Range("J19").Value = 446.75
$ws.Range("K19").Value = 1085.8572
$ws.Range("L19").Value = 446.75
$ws.Range("M19").Value = -910.8571999999999
$ws.Range("N19").Value = -796.75
$ws.Range("H41").Value = 433.66666
$ws.Range("J41").Value = 282.5
$ws.Range("L41").Value = 282.5
$ws.Range("N41").Value = -1162.5
$ws.Range("H100").Value = 4021.5
$ws.Range("I100").Value = 2036.8572
$ws.Range("K100").Value = 2036.8572
$ws.Range("M100").Value = -1495.8572
$ws.Range("H132").Value = 4671.75
$ws.Range("I132").Value = 3818.111
$ws.Range("K132").Value = 11454.333
$ws.Range("M132").Value = -8924.332999999999
$ws.Range("H137").Value = 1440.9615
$ws.Range("I137").Value = 1330.6364
$ws.Range("K137").Value = 3991.9092
$ws.Range("M137").Value = -1441.9092

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 2000000
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H61").Value = 1923.138
$ws.Range("I61").Value = 1658.5217
$ws.Range("J61").Value = 2937.5
$ws.Range("K61").Value = 1658.5217
$ws.Range("L61").Value = 2937.5
$ws.Range("M61").Value = -1446.5217
$ws.Range("N61").Value = -3361.5
$ws.Range("H74").Value = 1860.3392
$ws.Range("I74").Value = 1856.7858
$ws.Range("J74").Value = 1871
$ws.Range("K74").Value = 1856.7858
$ws.Range("L74").Value = 1871
$ws.Range("M74").Value = -982.7858000000001
$ws.Range("N74").Value = -3619
$ws.Range("H77").Value = 1860.3392
$ws.Range("I77").Value = 1856.7858
$ws.Range("J77").Value = 1871
$ws.Range("K77").Value = 9283.929
$ws.Range("L77").Value = 9355
$ws.Range("M77").Value = -4915.929
$ws.Range("N77").Value = -18091
$ws.Range("H97").Value = 822.1875
$ws.Range("I97").Value = 671.1
$ws.Range("K97").Value = 671.1
$ws.Range("M97").Value = -175.1
$ws.Range("H102").Value = 5343.2
$ws.Range("I102").Value = 2820.0833
$ws.Range("K102").Value = 2820.0833
$ws.Range("M102").Value = -1198.0833
$ws.Range("H122").Value = 2763.4
$ws.Range("I122").Value = 2678.7896
$ws.Range("K122").Value = 8036.3688
$ws.Range("M122").Value = -5586.3688
$ws.Range("H132").Value = 4673.756
$ws.Range("I132").Value = 4968.2163
$ws.Range("J132").Value = 1950
$ws.Range("K132").Value = 14904.6489
$ws.Range("L132").Value = 5850
$ws.Range("M132").Value = -12374.6489
$ws.Range("N132").Value = -10910
$ws.Range("H136").Value = 1923.138
$ws.Range("I136").Value = 1658.5217
$ws.Range("J136").Value = 2937.5
$ws.Range("K136").Value = 4975.5651
$ws.Range("L136").Value = 8812.5
$ws.Range("M136").Value = -2425.5651
$ws.Range("N136").Value = -13912.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3826
$ws.Range("I105").Value = 4586.4
$ws.Range("K105").Value = 4586.4
$ws.Range("M105").Value = -2839.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3691.182
$ws.Range("I31").Value = 1224.6666
$ws.Range("K31").Value = 1224.6666
$ws.Range("M31").Value = -929.6666
$ws.Range("H34").Value = 3691.182
$ws.Range("I34").Value = 1224.6666
$ws.Range("K34").Value = 1224.6666
$ws.Range("M34").Value = -1022.6666
$ws.Range("H58").Value = 1661.1666
$ws.Range("I58").Value = 1713
$ws.Range("J58").Value = 1588.6
$ws.Range("K58").Value = 1713
$ws.Range("L58").Value = 1588.6
$ws.Range("M58").Value = -1510
$ws.Range("N58").Value = -1994.6
$ws.Range("H99").Value = 1000
$ws.Range("J99").Value = 1000
$ws.Range("L99").Value = 1000
$ws.Range("N99").Value = -3996
$ws.Range("H103").Value = 9888.700000000001
$ws.Range("I103").Value = 9888.700000000001
$ws.Range("K103").Value = 9888.700000000001
$ws.Range("M103").Value = -8716.700000000001
$ws.Range("H126").Value = 1000
$ws.Range("J126").Value = 1000
$ws.Range("L126").Value = 3000
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 1818.4
$ws.Range("I132").Value = 1780.5862
$ws.Range("J132").Value = 2001.1666
$ws.Range("K132").Value = 5341.7586
$ws.Range("L132").Value = 6003.4998
$ws.Range("M132").Value = -2811.7586
$ws.Range("N132").Value = -11063.4998
$ws.Range("H134").Value = 1249.7833
$ws.Range("I134").Value = 1266.2559
$ws.Range("J134").Value = 1208.1177
$ws.Range("K134").Value = 3798.7677
$ws.Range("L134").Value = 3624.3531
$ws.Range("M134").Value = -1263.7677
$ws.Range("N134").Value = -8694.3531
$ws.Range("H136").Value = 1661.1666
$ws.Range("I136").Value = 1713
$ws.Range("J136").Value = 1588.6
$ws.Range("K136").Value = 5139
$ws.Range("L136").Value = 4765.799999999999
$ws.Range("M136").Value = -2589
$ws.Range("N136").Value = -9865.799999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 79.59
$ws.Range("I4").Value = 78.053764
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 234.161292
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -122.161292
$ws.Range("N4").Value = -524
$ws.Range("H59").Value = 998.3333
$ws.Range("I59").Value = 998.3333
$ws.Range("K59").Value = 2994.9999
$ws.Range("M59").Value = -2454.9999
$ws.Range("H117").Value = 2821.1667
$ws.Range("J117").Value = 1933.3334
$ws.Range("L117").Value = 5800.0002
$ws.Range("N117").Value = -12684.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3849.3333
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 4359.2
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 13077.6
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -17977.6
$ws.Range("H132").Value = 2280.925
$ws.Range("I132").Value = 2275.3076
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 6825.9228
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -4295.9228
$ws.Range("N132").Value = -12560
$ws.Range("H139").Value = 99492.336
$ws.Range("J139").Value = 99492.336
$ws.Range("L139").Value = 99492.336
$ws.Range("N139").Value = -109772.336

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1759.8182
$ws.Range("I55").Value = 539.8570999999999
$ws.Range("J55").Value = 2329.1333
$ws.Range("K55").Value = 539.8570999999999
$ws.Range("L55").Value = 2329.1333
$ws.Range("M55").Value = -366.8570999999999
$ws.Range("N55").Value = -2675.1333
$ws.Range("H76").Value = 17398
$ws.Range("J76").Value = 17398
$ws.Range("L76").Value = 17398
$ws.Range("N76").Value = -18074
$ws.Range("H79").Value = 17398
$ws.Range("J79").Value = 17398
$ws.Range("L79").Value = 17398
$ws.Range("N79").Value = -19738
$ws.Range("H93").Value = 5424.684
$ws.Range("I93").Value = 2739.25
$ws.Range("J93").Value = 7377.727
$ws.Range("K93").Value = 2739.25
$ws.Range("L93").Value = 7377.727
$ws.Range("M93").Value = -1491.25
$ws.Range("N93").Value = -9873.726999999999
$ws.Range("H132").Value = 5096.8335
$ws.Range("I132").Value = 4758.8335
$ws.Range("K132").Value = 14276.5005
$ws.Range("M132").Value = -11746.5005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 43285.57
$ws.Range("I19").Value = 60400
$ws.Range("J19").Value = 499.5
$ws.Range("K19").Value = 60400
$ws.Range("L19").Value = 499.5
$ws.Range("M19").Value = -60226
$ws.Range("N19").Value = -847.5
$ws.Range("H96").Value = 3818
$ws.Range("I96").Value = 2809
$ws.Range("K96").Value = 2809
$ws.Range("M96").Value = -1436
$ws.Range("H122").Value = 6762.4546
$ws.Range("I122").Value = 3477.4
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 10432.2
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -7982.200000000001
$ws.Range("N122").Value = -33400

Write-Output "applied changes"